# Apply the changes described in the commit "break out stock.yaml completed"
#
# 1) On the "week" sheet, cells D763:D776 currently hold the BSE code as
#    inline text; convert them to genuine numeric values (same digits).
# 2) On the "month" sheet, append 28 new data rows (53-80) describing the
#    latest month's top stocks, matching the existing table layout:
#    A=sr, B=nsecode, C=name, D=bsecode, E=per_chg, F=close, G=volume,
#    H=timeframe, I=Date Time

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Fix the D column type on the "week" sheet for rows 763-776
# ---------------------------------------------------------------
$weekSheet = $wb.Worksheets.Item("week")

$bseCodes = @{
    763 = 505200
    764 = 540115
    765 = 540762
    766 = 500520
    767 = 500257
    768 = 524715
    769 = 541450
    770 = 539254
    771 = 542066
    772 = 512070
    773 = 532898
    774 = 500477
    775 = 532234
    776 = 500183
}

foreach ($row in $bseCodes.Keys) {
    $weekSheet.Cells.Item($row, 4).Value = $bseCodes[$row]
}

# ---------------------------------------------------------------
# 2. Append the new rows to the "month" sheet
# ---------------------------------------------------------------
$monthSheet = $wb.Worksheets.Item("month")

$newRows = @(
    @(1,  "HAL",        "Hindustan Aeronautics Ltd",                  "541154", 2.38,  4178.35, 901281,    "month", "31/12/2024 21:37:47"),
    @(2,  "ANGELONE",   "Angel One Ltd",                              "",       0.1,   2931.65, 823118,    "month", "31/12/2024 21:37:47"),
    @(3,  "GODREJPROP", "Godrej Properties Limited",                  "533150", -1.38, 2786.5,  358989,    "month", "31/12/2024 21:37:47"),
    @(4,  "GRASIM",     "Grasim Industries Limited",                  "500300", -0.03, 2442.7,  530535,    "month", "31/12/2024 21:37:47"),
    @(5,  "INDIAMART",  "Indiamart Intermesh Ltd",                    "542726", -0.88, 2248.8,  270289,    "month", "31/12/2024 21:37:47"),
    @(6,  "ACC",        "Acc Limited",                                "500410", -0.36, 2052.5,  162386,    "month", "31/12/2024 21:37:47"),
    @(7,  "ASTRAL",     "Astral Poly Technik Limited",                "532830", -0.51, 1652.25, 424069,    "month", "31/12/2024 21:37:47"),
    @(8,  "KPITTECH",   "KPIT Technologies Ltd",                      "542651", -0.37, 1464.7,  615533,    "month", "31/12/2024 21:37:47"),
    @(9,  "BATAINDIA",  "Bata India Limited",                         "500043", -1.01, 1374.9,  193306,    "month", "31/12/2024 21:37:47"),
    @(10, "RELIANCE",   "Reliance Industries Limited",                "500325", 0.39,  1215.45, 6405475,   "month", "31/12/2024 21:37:47"),
    @(11, "AXISBANK",   "Axis Bank Limited",                          "532215", -0.49, 1064.7,  5292136,   "month", "31/12/2024 21:37:47"),
    @(12, "TATAMOTORS", "Tata Motors Limited",                        "500570", 0.89,  740.15,  7092699,   "month", "31/12/2024 21:37:47"),
    @(13, "JSL",        "Jindal Stainless Limited",                   "532508", 1.79,  699.2,   1157663,   "month", "31/12/2024 21:37:47"),
    @(14, "VBL",        "Varun Beverages Limited",                    "540180", -0.77, 638.5,   3444309,   "month", "31/12/2024 21:37:47"),
    @(15, "LICHSGFIN",  "Lic Housing Finance Limited",                "500253", 1.23,  598.05,  660447,    "month", "31/12/2024 21:37:47"),
    @(16, "AMBUJACEM",  "Ambuja Cements Limited",                     "500425", -1.62, 535.8,   1065733,   "month", "31/12/2024 21:37:47"),
    @(17, "RECLTD",     "Rural Electrification Corporation Limited",  "532955", 0.85,  500.7,   3906855,   "month", "31/12/2024 21:37:47"),
    @(18, "PFC",        "Power Finance Corporation Limited",          "532810", 1.83,  448.5,   6078911,   "month", "31/12/2024 21:37:47"),
    @(19, "BEL",        "Bharat Electronics Limited",                 "500049", 2.9,   293.15,  15725947,  "month", "31/12/2024 21:37:47"),
    @(20, "NCC",        "Ncc Limited",                                "500294", 1.33,  273.9,   2595902,   "month", "31/12/2024 21:37:47"),
    @(21, "HUDCO",      "Housing and Urban Development Corporation",  "540530", 4.04,  234.71,  8758607,   "month", "31/12/2024 21:37:47"),
    @(22, "BHEL",       "Bharat Heavy Electricals Limited",           "500103", 1.68,  229.4,   8401290,   "month", "31/12/2024 21:37:47"),
    @(23, "GAIL",       "Gail (india) Limited",                       "532155", 1.19,  190.98,  8397949,   "month", "31/12/2024 21:37:47"),
    @(24, "MANAPPURAM", "Manappuram Finance Limited",                 "531213", -0.33, 188.45,  6864819,   "month", "31/12/2024 21:37:47"),
    @(25, "IRFC",       "Indian Railway Finance Corporation Ltd",     "543257", -2.81, 149.04,  41308416,  "month", "31/12/2024 21:37:47"),
    @(26, "LTF",        "L&T Finance Ltd",                            "533519", -0.49, 135.63,  3059127,   "month", "31/12/2024 21:37:47"),
    @(27, "NHPC",       "Nhpc Limited",                               "533098", 1.18,  80.69,   24906446,  "month", "31/12/2024 21:37:47"),
    @(28, "GMRAIRPORT", "GMR Airports Ltd",                           "532754", 1.6,   78.56,   14038476,  "month", "31/12/2024 21:37:47")
)

$startRow = 53
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $monthSheet.Cells.Item($r, 1).Value = $data[0]
    $monthSheet.Cells.Item($r, 2).Value = $data[1]
    $monthSheet.Cells.Item($r, 3).Value = $data[2]
    if ($data[3] -eq "") {
        $monthSheet.Cells.Item($r, 4).Value = ""
    } else {
        # bsecode must stay as text (it is an inline string in the target
        # workbook), so prefix with an apostrophe to stop Excel from
        # auto-converting the numeric-looking text into a number.
        $monthSheet.Cells.Item($r, 4).Value = "'" + $data[3]
    }
    $monthSheet.Cells.Item($r, 5).Value = $data[4]
    $monthSheet.Cells.Item($r, 6).Value = $data[5]
    $monthSheet.Cells.Item($r, 7).Value = $data[6]
    $monthSheet.Cells.Item($r, 8).Value = $data[7]
    $monthSheet.Cells.Item($r, 9).Value = $data[8]
}
